$d = $word.ActiveDocument
$d.Paragraphs.Item(1).Range.Delete()
